# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" fund-holdings sheet (positioned between "总计"
# and the existing "2022-Q1" sheet) and records its summary row on the
# "总计" sheet, shifting the existing "2022-Q1"/"2021-Q2" summary rows
# down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a 2022-Q4 row right after the header,
#    pushing the existing 2022-Q1 / 2021-Q2 rows down one row each.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Snapshot the current row2 / row3 contents before we overwrite anything.
$oldRow2B = $total.Range("B2").Text
$oldRow2C = $total.Range("C2").Value2
$oldRow2D = $total.Range("D2").Value2

$oldRow3B = $total.Range("B3").Text
$oldRow3C = $total.Range("C3").Value2
$oldRow3D = $total.Range("D3").Value2

# Row 4 <- old row 3 (2021-Q2), keep the same bold/centered/bordered look
# that the other index cells (A2, A3) already use.
$total.Range("A4").Value = 2
$total.Range("B4").Value = $oldRow3B
$total.Range("C4").Value = $oldRow3C
$total.Range("D4").Value = $oldRow3D
$total.Range("A4").Font.Bold = $true
$total.Range("A4").HorizontalAlignment = -4108
$total.Range("A4").VerticalAlignment = -4160
$total.Range("A4").Borders.LineStyle = 1

# Row 3 <- old row 2 (2022-Q1).
$total.Range("A3").Value = 1
$total.Range("B3").Value = $oldRow2B
$total.Range("C3").Value = $oldRow2C
$total.Range("D3").Value = $oldRow2D

# Row 2 <- new 2022-Q4 summary.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with the fund holdings detail, inserted
#    immediately before the existing "2022-Q1" sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q4 = $wb.Worksheets.Add($q1)
$q4.Name = "2022-Q4"

# Match the page-setup margins used by the sibling fund-detail sheets.
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $q4.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}
$hdr = $q4.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$rows = @(
    @("009514", "创金合信同顺创业板精选股票C", "0.15", "91.10", "1.54", "0.0023", 8),
    @("009513", "创金合信同顺创业板精选股票A", "0.09", "91.10", "1.54", "0.0014", 8)
)

$r = 2
foreach ($row in $rows) {
    $idxCell = $q4.Range("A" + $r)
    $idxCell.Value = $r - 2
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    # Fund code keeps its leading zeros (e.g. "009514") -- must stay text.
    $q4.Range("B" + $r).NumberFormat = "@"
    $q4.Range("B" + $r).Value = $row[0]
    $q4.Range("C" + $r).Value = $row[1]

    # Columns D-G keep their original text representation (mixed
    # percentages / decimals) instead of being auto-coerced to numbers.
    $q4.Range("D" + $r).NumberFormat = "@"
    $q4.Range("D" + $r).Value = $row[2]
    $q4.Range("E" + $r).NumberFormat = "@"
    $q4.Range("E" + $r).Value = $row[3]
    $q4.Range("F" + $r).NumberFormat = "@"
    $q4.Range("F" + $r).Value = $row[4]
    $q4.Range("G" + $r).NumberFormat = "@"
    $q4.Range("G" + $r).Value = $row[5]

    $q4.Range("H" + $r).Value = $row[6]

    $r = $r + 1
}

# Keep "总计" as the active sheet, same as before the edit.
$total.Activate()
